$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; existing D:K data shifts to E:L
$ws.Columns("D:D").Insert()

# Copy number/date formatting from the (now-shifted) old column D, sitting in E, back into new D
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate new column D with the new quarter (2018-09-30) figures
$ws.Cells.Item(7, 4).Value = 43373
$ws.Cells.Item(8, 4).Value = 9000
$ws.Cells.Item(9, 4).Value = 5700
$ws.Cells.Item(10, 4).Value = 3300
$ws.Cells.Item(12, 4).Value = 400
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(17, 4).Value = 8600
$ws.Cells.Item(18, 4).Value = 400
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(21, 4).Value = 500
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(23, 4).Value = 300
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 4).Value = 300
$ws.Cells.Item(27, 4).Value = 300
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(33, 4).Value = 300
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(35, 4).Value = 300
$ws.Cells.Item(38, 4).Value = 43373
$ws.Cells.Item(41, 4).Value = 200
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(43, 4).Value = 4400
$ws.Cells.Item(44, 4).Value = 6300
$ws.Cells.Item(45, 4).Value = 700
$ws.Cells.Item(46, 4).Value = 11600
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(48, 4).Value = 300
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(52, 4).Value = 200
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(54, 4).Value = 12100
$ws.Cells.Item(57, 4).Value = 3500
$ws.Cells.Item(58, 4).Value = 1900
$ws.Cells.Item(59, 4).Value = 2100
$ws.Cells.Item(60, 4).Value = 7600
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(66, 4).Value = 7600
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(72, 4).Value = -36500
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(76, 4).Value = 4500
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(80, 4).Value = 43373
$ws.Cells.Item(81, 4).Value = 300
$ws.Cells.Item(83, 4).Value = 100
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(89, 4).Value = -1500
$ws.Cells.Item(91, 4).Value = -200
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(94, 4).Value = -200
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(100, 4).Value = 1800
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(102, 4).Value = 100
